{"js": "// Update the worksheet date and all 25 division problems to the new values.\nconst replacements = [\n  [\"2024-05-19 Sunday\", \"2024-05-20 Monday\"],\n  [\"757\u00f72=\", \"820\u00f73=\"],\n  [\"255\u00f73=\", \"713\u00f75=\"],\n  [\"678\u00f73=\", \"586\u00f77=\"],\n  [\"297\u00f75=\", \"826\u00f79=\"],\n  [\"493\u00f75=\", \"526\u00f72=\"],\n  [\"368\u00f78=\", \"927\u00f74=\"],\n  [\"246\u00f75=\", \"444\u00f78=\"],\n  [\"835\u00f74=\", \"129\u00f74=\"],\n  [\"791\u00f78=\", \"228\u00f77=\"],\n  [\"466\u00f74=\", \"806\u00f72=\"],\n  [\"670\u00f74=\", \"698\u00f74=\"],\n  [\"579\u00f77=\", \"888\u00f74=\"],\n  [\"432\u00f77=\", \"688\u00f76=\"],\n  [\"274\u00f76=\", \"126\u00f76=\"],\n  [\"337\u00f75=\", \"454\u00f77=\"],\n  [\"271\u00f77=\", \"663\u00f75=\"],\n  [\"888\u00f77=\", \"845\u00f72=\"],\n  [\"521\u00f72=\", \"127\u00f78=\"],\n  [\"643\u00f73=\", \"249\u00f75=\"],\n  [\"785\u00f79=\", \"999\u00f76=\"],\n  [\"818\u00f74=\", \"228\u00f76=\"],\n  [\"716\u00f78=\", \"836\u00f75=\"],\n  [\"218\u00f75=\", \"317\u00f74=\"],\n  [\"192\u00f73=\", \"270\u00f75=\"],\n  [\"242\u00f78=\", \"185\u00f72=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Search text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and all 25 division problems to the new values.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-05-19 Sunday\", \"2024-05-20 Monday\"),\n    @(\"757\u00f72=\", \"820\u00f73=\"),\n    @(\"255\u00f73=\", \"713\u00f75=\"),\n    @(\"678\u00f73=\", \"586\u00f77=\"),\n    @(\"297\u00f75=\", \"826\u00f79=\"),\n    @(\"493\u00f75=\", \"526\u00f72=\"),\n    @(\"368\u00f78=\", \"927\u00f74=\"),\n    @(\"246\u00f75=\", \"444\u00f78=\"),\n    @(\"835\u00f74=\", \"129\u00f74=\"),\n    @(\"791\u00f78=\", \"228\u00f77=\"),\n    @(\"466\u00f74=\", \"806\u00f72=\"),\n    @(\"670\u00f74=\", \"698\u00f74=\"),\n    @(\"579\u00f77=\", \"888\u00f74=\"),\n    @(\"432\u00f77=\", \"688\u00f76=\"),\n    @(\"274\u00f76=\", \"126\u00f76=\"),\n    @(\"337\u00f75=\", \"454\u00f77=\"),\n    @(\"271\u00f77=\", \"663\u00f75=\"),\n    @(\"888\u00f77=\", \"845\u00f72=\"),\n    @(\"521\u00f72=\", \"127\u00f78=\"),\n    @(\"643\u00f73=\", \"249\u00f75=\"),\n    @(\"785\u00f79=\", \"999\u00f76=\"),\n    @(\"818\u00f74=\", \"228\u00f76=\"),\n    @(\"716\u00f78=\", \"836\u00f75=\"),\n    @(\"218\u00f75=\", \"317\u00f74=\"),\n    @(\"192\u00f73=\", \"270\u00f75=\"),\n    @(\"242\u00f78=\", \"185\u00f72=\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $found = $find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n    if (-not $found) {\n        Write-Host \"Not found: $old\"\n    }\n}\n"}
